$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 54: becomes the closing row of the um2502 script block ---
# Copy the formatting from row 4 (same two-row-block "closing" style: s=6/7)
$ws.Range("A4:E4").Copy()
$ws.Range("A54:E54").PasteSpecial(-4122)

# --- Row 55: brand-new row, start of the us0101 script block (style 4/5) ---
$ws.Range("A53:E53").Copy()
$ws.Range("A55:E55").PasteSpecial(-4122)

# Fill the new cell values in the same order the shared-string table grew in
# the target workbook (A54, C55, A55, D55, E55 -- B55 is a plain number).
$ws.Range("A54").Value = "SCRIPT/G01P03A/um2502.ssb"
$ws.Range("C55").Value = " Good luck on your graduation\nexam!"
$ws.Range("A55").Value = "SCRIPT/G01P03A/us0101.ssb"
$ws.Range("D55").Value = " Удачи вам на выпускном\nэкзамене!"
$ws.Range("E55").Value = " Ôäàœé âàí îà âúðôòëîïí\nüëèàíåîå!"
$ws.Range("B55").Value = 276

$ws.Rows.Item(54).RowHeight = 43.2
$ws.Rows.Item(55).RowHeight = 46.8

# --- View state: keep selection/scroll close to what the author left behind ---
$ws.Range("D59").Select()
$excel.ActiveWindow.ScrollRow = 53

Write-Host "done"
